# HelloWorld/Template.xlsx — flesh out the FileSheet "Paths" section and
# add a settings import from a serialized xml (per commit message):
#   - The "Paths" block gains an "Inflation" row, mirroring the
#     highlighted-label style already used for "Template" (D4).
#   - The old placeholder text that lived in E4 (a leftover absolute path
#     to this very template file) is cleared.
#   - The active selection on the File_Template sheet moves from D16 to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("File_Template")

# Give D5 the same "highlighted label" formatting as D4 (fill/border used
# for the Paths column labels) by copying D4's format onto it.
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Move the "Inflation" label from E4 (where it was a stray value) to the
# newly formatted D5 cell.
$ws.Range("E4").ClearContents()
$ws.Range("D5").Value = "Inflation"

# Update the sheet's active selection to D5.
$ws.Activate() | Out-Null
$ws.Range("D5").Select() | Out-Null
